$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 233, shifting existing rows 233:311 down to 234:312
$ws.Rows.Item(233).Insert()

# Fill in the new row 233 with the new weekly price entry
$ws.Cells.Item(233, 1).Value = 10
$ws.Cells.Item(233, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(233, 3).Value = "La Araucanía"
$ws.Cells.Item(233, 4).Value = 45120
$ws.Cells.Item(233, 5).Value = 9
$ws.Cells.Item(233, 6).Value = 100112013
$ws.Cells.Item(233, 7).Value = "Alcachofa"
$ws.Cells.Item(233, 8).Value = "Madrigal"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 300
$ws.Cells.Item(233, 11).Value = 16000
$ws.Cells.Item(233, 12).Value = 16000
$ws.Cells.Item(233, 13).Value = 16000
$ws.Cells.Item(233, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(233, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(233, 16).Value = 533
$ws.Cells.Item(233, 17).Value = 30
$ws.Cells.Item(233, 18).Value = "Hortaliza"
